$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header E1 from "Location" to "Country"
$ws.Range("E1").Value = "Country"

# Add two new lead rows (row 5 and row 6)
$ws.Range("A5").Value = "John Doe"
$ws.Range("B5").Value = 1234567890
$ws.Range("C5").Value = "Example Corp"
$ws.Range("D5").Value = "Product Company"
$ws.Range("E5").Value = "India"
$ws.Range("F5").Value = "Chennai"
$ws.Range("G5").Value = "One Time"
$ws.Range("H5").Value = "Google Adds, Marketing"
$ws.Range("I5").Value = 50000
$ws.Range("J5").Value = 45000
$ws.Range("K5").Value = 45443
$ws.Range("K5").NumberFormat = "yyyy-mm-dd"
$ws.Range("L5").Value = "High"
$ws.Range("M5").Value = "johndoe@example.com"
$ws.Range("N5").Value = "Fresh"
$ws.Range("O5").Value = "Initial contact made"
$ws.Range("P5").Value = "Follow up next week"
$ws.Range("Q5").Value = 45443
$ws.Range("Q5").NumberFormat = "yyyy-mm-dd"
$ws.Range("R5").Value = $false
$ws.Range("S5").Value = 45443
$ws.Range("S5").NumberFormat = "yyyy-mm-dd"
$ws.Range("T5").Value = "admin"
$ws.Range("U5").Value = "admin"
$ws.Range("V5").Value = 45443
$ws.Range("V5").NumberFormat = "yyyy-mm-dd"

$ws.Range("A6").Value = "John Doe"
$ws.Range("B6").Value = 1234567890
$ws.Range("C6").Value = "Example Corp"
$ws.Range("D6").Value = "Product Company"
$ws.Range("E6").Value = "India"
$ws.Range("F6").Value = "Chennai"
$ws.Range("G6").Value = "One Time"
$ws.Range("H6").Value = "SMM"
$ws.Range("I6").Value = 50000
$ws.Range("J6").Value = 45000
$ws.Range("K6").Value = 45443
$ws.Range("K6").NumberFormat = "yyyy-mm-dd"
$ws.Range("L6").Value = "High"
$ws.Range("M6").Value = "johndoe@example.com"
$ws.Range("N6").Value = "Fresh"
$ws.Range("O6").Value = "Initial contact made"
$ws.Range("P6").Value = "Follow up next week"
$ws.Range("Q6").Value = 45443
$ws.Range("Q6").NumberFormat = "yyyy-mm-dd"
$ws.Range("R6").Value = $false
$ws.Range("S6").Value = 45443
$ws.Range("S6").NumberFormat = "yyyy-mm-dd"
$ws.Range("T6").Value = "admin"
$ws.Range("U6").Value = "admin"
$ws.Range("V6").Value = 45443
$ws.Range("V6").NumberFormat = "yyyy-mm-dd"
